$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: "Sending cluster" changes from Resolving-Mac to MuSCs for all data rows ---
$ws.Range("A2").Value2 = "MuSCs"
$ws.Range("A3").Value2 = "MuSCs"
$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("A5").Value2 = "MuSCs"

# --- Column D: "Target cluster" - rows 4 and 5 swap (MuSCs <-> Resolving-Mac) ---
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("D5").Value2 = "Resolving-Mac"

# --- Row 2 (Target cluster = ECs) updated TPM-derived metrics ---
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.06815433333333333
$ws.Range("H2").Value2 = 0.204463
$ws.Range("M2").Value2 = 13.17295566666667
$ws.Range("N2").Value2 = 39.518867
$ws.Range("O2").Value2 = 0.133784132206724
$ws.Range("P2").Value2 = 0.133784132206724
$ws.Range("Q2").Value2 = 0.8977940114912222
$ws.Range("R2").Value2 = 8.080146103421001
$ws.Range("S2").Value2 = 0.133784132206724
$ws.Range("T2").Value2 = 0.133784132206724

# --- Row 3 (Target cluster = FAPs) updated TPM-derived metrics ---
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.06815433333333333
$ws.Range("H3").Value2 = 0.204463
$ws.Range("O3").Value2 = 0.4382627974978752
$ws.Range("P3").Value2 = 0.4382627974978752
$ws.Range("Q3").Value2 = 2.941079099313444
$ws.Range("R3").Value2 = 26.469711893821
$ws.Range("S3").Value2 = 0.4382627974978752
$ws.Range("T3").Value2 = 0.4382627974978752

# --- Row 4 (Target cluster = MuSCs) updated TPM-derived metrics ---
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.06815433333333333
$ws.Range("H4").Value2 = 0.204463
$ws.Range("M4").Value2 = 21.06166566666667
$ws.Range("N4").Value2 = 63.184997
$ws.Range("O4").Value2 = 0.2139016281041017
$ws.Range("P4").Value2 = 0.2139016281041017
$ws.Range("Q4").Value2 = 1.435443782401222
$ws.Range("R4").Value2 = 12.918994041611
$ws.Range("S4").Value2 = 0.2139016281041017
$ws.Range("T4").Value2 = 0.2139016281041017

# --- Row 5 (Target cluster = Resolving-Mac) updated TPM-derived metrics ---
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.06815433333333333
$ws.Range("H5").Value2 = 0.204463
$ws.Range("M5").Value2 = 21.076417
$ws.Range("N5").Value2 = 63.229251
$ws.Range("O5").Value2 = 0.214051442191299
$ws.Range("P5").Value2 = 0.214051442191299
$ws.Range("Q5").Value2 = 1.436449149690334
$ws.Range("R5").Value2 = 12.928042347213
$ws.Range("S5").Value2 = 0.214051442191299
$ws.Range("T5").Value2 = 0.214051442191299
